# -------------------------------------------------------------------------
# Tyres worksheet refresh: tyre catalogue re-keyed + new sale date
# (recovered after second total system crash, 14/11/2023 -> serial 45244)
# -------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13: full records (A,B,C + E,F,G,H,I,J mirrored/lookup columns)
# Columns: A=Model, B=Size, C=Param(raw), E=Size, F=Model, G=Param(list), H=Sales, I=Date, J=Contragent
$fullRows = @(
    [pscustomobject]@{ Row=2; A="ФБел-283"; B="35/65-33"; C="42 груз сер"; E="35/65-33"; F="ФБел-283"; G="42, груз, сер"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=3; A="Бел-122"; B="24.00R35"; C="груз C сер H Type LS-2"; E="205/55R16"; F="BEL-262"; G="б/к, сер, легк"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=4; A="Бел-202"; B="24.00R35"; C="210B C сер H Type LS-2"; E="205/55R16"; F="BEL-317"; G="б/к, сер, легк"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=5; A="Бел-212"; B="24.00R35"; C="груз C сер H Type LS-2"; E="205/55R16"; F="BEL-317S"; G="сер, ошип"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=6; A="Бел-200"; B="21.00R35"; C="202B C сер H Type LS-2"; E="24.00R35"; F="Бел-202"; G="210B, C, сер, H, Type, LS-2"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=7; A="Бел-210"; B="21.00R35"; C="202B C сер H Type LS-2"; E="24.00R35"; F="Бел-212"; G="груз, C, сер, H, Type, LS-2"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=8; A="BEL-248"; B="14.00R20"; C="б/к груз сер"; E="21.00R35"; F="Бел-200"; G="202B, C, сер, H, Type, LS-2"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=9; A="BEL-288"; B="12.00R20"; C="груз сер"; E="21.00R35"; F="Бел-210"; G="202B, C, сер, H, Type, LS-2"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=10; A="Бел-95"; B="16.00R20"; C="173G б/к груз сер"; E="14.00R20"; F="BEL-248"; G="б/к, груз, сер"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=11; A="BEL-405"; B="395/85R20"; C="168J груз сер"; E="175/70R13"; F="Бел-103"; G="б/к, сер, легк"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=12; A="Бел-145"; B="445/65R22.5"; C="б/к груз сер"; E="175/70R13"; F="Бел-100"; G="б/к, сер, легк"; H=2; I=45244; J="нет данных" },
    [pscustomobject]@{ Row=13; A="Бел-230"; B="355/65-15"; C="сер"; E="195/65R15"; F="Бел-119"; G="сер, легк"; H=2; I=45244; J="нет данных" }
)

foreach ($r in $fullRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A   # A
    $ws.Cells.Item($r.Row, 2).Value = $r.B   # B
    $ws.Cells.Item($r.Row, 3).Value = $r.C   # C
    $ws.Cells.Item($r.Row, 5).Value = $r.E   # E
    $ws.Cells.Item($r.Row, 6).Value = $r.F   # F
    $ws.Cells.Item($r.Row, 7).Value = $r.G   # G
    $ws.Cells.Item($r.Row, 8).Value = $r.H   # H
    $ws.Cells.Item($r.Row, 9).Value = $r.I   # I (date serial; cell already carries the date style)
    $ws.Cells.Item($r.Row, 10).Value = $r.J  # J
}

# Rows 14-22: only A,B,C retained; the old E:J lookup data for these rows is gone
$shortRows = @(
    [pscustomobject]@{ Row=14; A="Бел-230М"; B="355/65-15"; C="сер" },
    [pscustomobject]@{ Row=15; A="BEL-262"; B="205/55R16"; C="б/к сер легк" },
    [pscustomobject]@{ Row=16; A="BEL-317"; B="205/55R16"; C="б/к сер легк" },
    [pscustomobject]@{ Row=17; A="BEL-317S"; B="205/55R16"; C="сер ошип" },
    [pscustomobject]@{ Row=18; A="BEL-277"; B="205/60R16"; C="б/к сер легк" },
    [pscustomobject]@{ Row=19; A="Бел-103"; B="175/70R13"; C="б/к сер легк" },
    [pscustomobject]@{ Row=20; A="Бел-100"; B="175/70R13"; C="б/к сер легк" },
    [pscustomobject]@{ Row=21; A="Ф-35-1"; B="11.2-20"; C="8 168J сх сер" },
    [pscustomobject]@{ Row=22; A="Бел-119"; B="195/65R15"; C="сер легк" }
)

foreach ($r in $shortRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Range($ws.Cells.Item($r.Row, 4), $ws.Cells.Item($r.Row, 10)).Clear()
}

# Row 23 (superDuper / 155/65R13 / легк) was dropped entirely from the catalogue
$ws.Range("A23:J23").Clear()
